# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the "8f18abb7-81b6-4f07-8c7b-0431a06db6b0" record, in
# both the zh-cn and de-de language sheets (and, by virtue of the shared
# "Status" string, the Overview roll-up sheet as well).

$wb = $excel.ActiveWorkbook

$zhCn = $wb.Worksheets.Item("zh-cn")
$deDe = $wb.Worksheets.Item("de-de")

# The Status column value "Ready for handoff" is shared across the
# Overview, zh-cn and de-de sheets for this record, so replace it
# workbook-wide in one go rather than sheet-by-sheet.
foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Ready for handoff", "Handback transform failed")
}

# Row 3 is the "8f18abb7-81b6-4f07-8c7b-0431a06db6b0" record in each
# language table. Column K = Error Detail - fill in why the handback
# transform failed.

$zhCn.Range("K3").Value = "Handback file name: gecovcql.5yf is different with handoff file name: 8f18abb7-81b6-4f07-8c7b-0431a06db6b0.7fb4a54f3c3b9f8054cd75b381d4e39399c9264c.zh-cn."

$deDe.Range("K3").Value = "Handback file name: gecovcql.5yf is different with handoff file name: 8f18abb7-81b6-4f07-8c7b-0431a06db6b0.7fb4a54f3c3b9f8054cd75b381d4e39399c9264c.de-de."
